# Fruta / hortaliza, semanal
# Insert two new price records at the top of the Ciruela (plum) data block
# (rows 236-237), shifting the existing rows 236-275 down to 238-277.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 236 (old row 236 -> new row 238, etc.)
$ws.Rows(236).Insert()
$ws.Rows(237).Insert()

# --- New row 236: Black Amber ---
$ws.Range("A236").Value = 10
$ws.Range("B236").Value = "Vega Modelo de Temuco"
$ws.Range("C236").Value = "La Araucanía"
$ws.Range("D236").Value = 44946
$ws.Range("E236").Value = 9
$ws.Range("F236").Value = "Fruta"
$ws.Range("G236").Value = 100103
$ws.Range("H236").Value = "Frutos de hueso (carozo)"
$ws.Range("I236").Value = 100103002
$ws.Range("J236").Value = "Ciruela"
$ws.Range("K236").Value = "Black Amber"
$ws.Range("L236").Value = "Primera"
$ws.Range("M236").Value = 125
$ws.Range("N236").Value = 18000
$ws.Range("O236").Value = 18000
$ws.Range("P236").Value = 18000
$ws.Range("Q236").Value = "$/bandeja 18 kilos granel"
$ws.Range("R236").Value = "Región de O'Higgins"
$ws.Range("S236").Value = 1000
$ws.Range("T236").Value = 18

# --- New row 237: Pink Delight ---
$ws.Range("A237").Value = 10
$ws.Range("B237").Value = "Vega Modelo de Temuco"
$ws.Range("C237").Value = "La Araucanía"
$ws.Range("D237").Value = 44946
$ws.Range("E237").Value = 9
$ws.Range("F237").Value = "Fruta"
$ws.Range("G237").Value = 100103
$ws.Range("H237").Value = "Frutos de hueso (carozo)"
$ws.Range("I237").Value = 100103002
$ws.Range("J237").Value = "Ciruela"
$ws.Range("K237").Value = "Pink Delight"
$ws.Range("L237").Value = "Primera"
$ws.Range("M237").Value = 115
$ws.Range("N237").Value = 15000
$ws.Range("O237").Value = 15000
$ws.Range("P237").Value = 15000
$ws.Range("Q237").Value = "$/bandeja 18 kilos granel"
$ws.Range("R237").Value = "Región de O'Higgins"
$ws.Range("S237").Value = 833
$ws.Range("T237").Value = 18
